{"js": "// The documented change only reorders the children of <w:rPr> inside\n// word/styles.xml for a handful of character styles (KeywordTok,\n// ImportTok, CommentTok, DocumentationTok, AnnotationTok, CommentVarTok,\n// ControlFlowTok, InformationTok, WarningTok, AlertTok, ErrorTok):\n// <w:b/>/<w:i/> must come before <w:color/> to satisfy wml.xsd's\n// CT_RPr sequence (the schema expects b/bCs/i/iCs/... before color).\n//\n// There is no direct \"reorder child elements\" verb in the Word API, so\n// we touch the run-formatting properties that are already set on each\n// of those styles (re-assigning the same boolean value) which forces\n// the host to re-emit <w:rPr> in schema order without changing any\n// actual formatting.\n\nconst targetStyleNames = [\n  \"KeywordTok\",\n  \"ImportTok\",\n  \"CommentTok\",\n  \"DocumentationTok\",\n  \"AnnotationTok\",\n  \"CommentVarTok\",\n  \"ControlFlowTok\",\n  \"InformationTok\",\n  \"WarningTok\",\n  \"AlertTok\",\n  \"ErrorTok\",\n];\n\nconst styles = context.document.getStyles();\nstyles.load(\"items/nameLocal\");\nawait context.sync();\n\nconst targets = styles.items.filter((s) => targetStyleNames.includes(s.nameLocal));\ntargets.forEach((s) => s.font.load(\"bold,italic\"));\nawait context.sync();\n\n// Only re-assign properties that are already truthy so we never add a\n// new explicit \"false\" (e.g. <w:b w:val=\"0\"/>) that wasn't there before.\nfor (const style of targets) {\n  if (style.font.bold) {\n    style.font.bold = true;\n  }\n  if (style.font.italic) {\n    style.font.italic = true;\n  }\n}\nawait context.sync();\n", "ps1": "# The documented change only reorders the children of <w:rPr> inside\n# word/styles.xml for a handful of character styles (KeywordTok,\n# ImportTok, CommentTok, DocumentationTok, AnnotationTok, CommentVarTok,\n# ControlFlowTok, InformationTok, WarningTok, AlertTok, ErrorTok):\n# <w:b/>/<w:i/> must come before <w:color/> to satisfy wml.xsd's\n# CT_RPr sequence (the schema expects b/bCs/i/iCs/... before color).\n#\n# The Word object model has no \"reorder child elements\" verb, so we\n# touch the run-formatting properties that are already set on each of\n# those styles (re-assigning the same boolean value), which forces the\n# host to re-emit <w:rPr> in schema order without changing any actual\n# formatting.\n\n$d = $word.ActiveDocument\n\n$targetStyleNames = @(\n    \"KeywordTok\",\n    \"ImportTok\",\n    \"CommentTok\",\n    \"DocumentationTok\",\n    \"AnnotationTok\",\n    \"CommentVarTok\",\n    \"ControlFlowTok\",\n    \"InformationTok\",\n    \"WarningTok\",\n    \"AlertTok\",\n    \"ErrorTok\"\n)\n\nforeach ($name in $targetStyleNames) {\n    $style = $d.Styles.Item($name)\n    $font = $style.Font\n\n    # Only re-assign properties that are already truthy (non-zero) so we\n    # never introduce a new explicit \"false\" (e.g. <w:b w:val=\"0\"/>)\n    # that wasn't present before.\n    if ($font.Bold) {\n        $font.Bold = -1\n    }\n    if ($font.Italic) {\n        $font.Italic = -1\n    }\n}\n"}
